# Applies the "Delimiters slightly reworked + separate scripts done" edit.
#
# Summary of the change:
#  - input_k_constants_log10: header text changed, and the numeric log10(K)
#    values are now stored as text (string) cells instead of numbers.
#  - input_concentrations: the "eq"/"tot" header labels on B1/C1 are swapped,
#    and the numeric concentration values are now stored as text cells
#    (with new numeric values).
#  - equilibrium_concentrations / L_fractions / percent_error: numeric
#    values were recalculated (new numbers), column headers stay the same
#    text as before.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$val) {
    # Force Excel to store the value as a text (string) cell, even when the
    # text looks like a number, without leaving a lingering "Text" number
    # format / quote-prefix style on the cell.
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet: input_k_constants_log10
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("input_k_constants_log10")

$ws.Range("A1").Value = "log10.K."

Set-TextValue $ws.Range("A2") "3.14"
Set-TextValue $ws.Range("A3") "1.45"
Set-TextValue $ws.Range("A4") "2.79"
Set-TextValue $ws.Range("A5") "-8.9"
Set-TextValue $ws.Range("A6") "-13.88"

# ---------------------------------------------------------------------
# Sheet: input_concentrations
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("input_concentrations")

$ws.Range("B1").Value = "eq"
$ws.Range("C1").Value = "tot"

Set-TextValue $ws.Range("A3") "0.1"
Set-TextValue $ws.Range("B3") "0.001"
Set-TextValue $ws.Range("C3") "0.005"

Set-TextValue $ws.Range("A4") "0.01"
Set-TextValue $ws.Range("B4") "0.001"
Set-TextValue $ws.Range("C4") "0.005"

Set-TextValue $ws.Range("A5") "1e-04"
Set-TextValue $ws.Range("B5") "0.01"
Set-TextValue $ws.Range("C5") "0.005"

# ---------------------------------------------------------------------
# Sheet: equilibrium_concentrations
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("equilibrium_concentrations")

$ws.Range("A2").Value = 0.1
$ws.Range("B2").Value = 0.001
$ws.Range("C2").Value = 0.004860029427680299
$ws.Range("D2").Value = 0.1380384264602885
$ws.Range("E2").Value = 0.0001369742398441712
$ws.Range("F2").Value = [double]"2.996669854007388e-06"
$ws.Range("G2").Value = [double]"6.118414548574175e-11"
$ws.Range("H2").Value = [double]"1.318256738556401e-13"

$ws.Range("A3").Value = 0.01
$ws.Range("B3").Value = 0.001
$ws.Range("C3").Value = 0.004860028892444237
$ws.Range("D3").Value = 0.01380384264602885
$ws.Range("E3").Value = 0.0001369742247591694
$ws.Range("F3").Value = [double]"2.996669523983508e-06"
$ws.Range("G3").Value = [double]"6.118413874751902e-10"
$ws.Range("H3").Value = [double]"1.318256738556402e-12"

$ws.Range("B4").Value = 0.01
$ws.Range("C4").Value = 0.003721593866991905
$ws.Range("D4").Value = 0.001380384264602885
$ws.Range("E4").Value = 0.001048887663182846
$ws.Range("F4").Value = 0.0002294716177345554
$ws.Range("G4").Value = [double]"4.68520909153343e-08"

# ---------------------------------------------------------------------
# Sheet: L_fractions
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("L_fractions")

$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 13803.84264602885
$ws.Range("E2").Value = 13.69742398441712
$ws.Range("F2").Value = 0.5993339708014775

$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 1380.384264602885
$ws.Range("E3").Value = 13.69742247591694
$ws.Range("F3").Value = 0.5993339047967016

$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 13.80384264602885
$ws.Range("E4").Value = 10.48887663182846
$ws.Range("F4").Value = 4.589432354691109

# ---------------------------------------------------------------------
# Sheet: percent_error
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("percent_error")

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = [double]"3.985626230093575e-10"

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = [double]"3.985687778082503e-10"

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = [double]"2.211772431870429e-16"
